# order add area flag 1=大连 0=沈阳
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# H37 was left blank while D37:G37 already carry the "check mark" (同上)
# string used throughout the sheet; tick the box to match the rest of the row.
$ws.Range("H37").Value = [char]0x221A
$ws.Range("G37").Copy()
$ws.Range("H37").PasteSpecial(-4122)  # xlPasteFormats - keep the centered √ style
$excel.CutCopyMode = $false

# Every sample login's password column was a run of sequential dummy numbers
# (123457, 123458, ...); reset them all to the same placeholder value 123456.
$ws.Range("F46:F57").Value = 123456

# Refresh the saved view: zoom out a bit and leave the window scrolled/selected
# further down the list (row ~30 visible, H50 selected) like the author left it.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 30
$ws.Range("H50").Select() | Out-Null
